$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantities for existing rows
$ws.Range("A2").Value = 26
$ws.Range("A3").Value = 30
$ws.Range("A4").Value = 3

# Add new inventory rows
$ws.Range("A5").Value = "utp patch kábel"
$ws.Range("C5").Value = "500m"

$ws.Range("A6").Value = "uplink kábel"
$ws.Range("C6").Value = "22m"

$ws.Range("A7").Value = "ether ch kábel"
$ws.Range("C7").Value = "8m"

$ws.Range("A8").Value = "router kábel"
$ws.Range("C8").Value = "4m"

$ws.Range("A9").Value = "WAN kábel"
$ws.Range("C9").Value = "4m"

$ws.Range("C9").Select()
